{"js": "// The document is a daily \"two-digit \u00f7 one-digit\" division drill sheet: a\n// centered date heading followed by a 5-column table whose populated rows\n// hold \"<dividend>\u00f7<divisor>=<quotient>, <remainder>\" answers. This commit\n// bumps the date by one day and swaps in a freshly generated set of answers\n// for every populated cell.\n//\n// Every old string below occurs exactly once in the document body, so for\n// each pair we search for the old text and replace the matched range's text\n// with the new value in place. insertText(..., \"Replace\") rewrites only the\n// matched run's text and keeps that run's existing formatting (font/size),\n// so no rPr/pPr needs to be touched.\nconst replacements = [\n  [\"2024-06-02 Sunday\", \"2024-06-03 Monday\"],\n  [\"20\u00f77=2, 6\", \"93\u00f75=18, 3\"],\n  [\"95\u00f78=11, 7\", \"73\u00f77=10, 3\"],\n  [\"46\u00f74=11, 2\", \"59\u00f78=7, 3\"],\n  [\"81\u00f77=11, 4\", \"78\u00f78=9, 6\"],\n  [\"68\u00f76=11, 2\", \"72\u00f78=9, 0\"],\n  [\"87\u00f79=9, 6\", \"15\u00f78=1, 7\"],\n  [\"48\u00f72=24, 0\", \"91\u00f77=13, 0\"],\n  [\"55\u00f75=11, 0\", \"25\u00f72=12, 1\"],\n  [\"69\u00f72=34, 1\", \"23\u00f72=11, 1\"],\n  [\"33\u00f74=8, 1\", \"75\u00f72=37, 1\"],\n  [\"13\u00f72=6, 1\", \"12\u00f72=6, 0\"],\n  [\"53\u00f74=13, 1\", \"24\u00f78=3, 0\"],\n  [\"25\u00f74=6, 1\", \"59\u00f73=19, 2\"],\n  [\"35\u00f79=3, 8\", \"36\u00f79=4, 0\"],\n  [\"54\u00f75=10, 4\", \"37\u00f78=4, 5\"],\n  [\"18\u00f75=3, 3\", \"90\u00f79=10, 0\"],\n  [\"30\u00f73=10, 0\", \"96\u00f74=24, 0\"],\n  [\"70\u00f74=17, 2\", \"96\u00f75=19, 1\"],\n  [\"20\u00f73=6, 2\", \"43\u00f78=5, 3\"],\n  [\"39\u00f76=6, 3\", \"55\u00f76=9, 1\"],\n  [\"78\u00f77=11, 1\", \"93\u00f77=13, 2\"],\n  [\"84\u00f75=16, 4\", \"76\u00f74=19, 0\"],\n  [\"13\u00f76=2, 1\", \"93\u00f75=18, 3\"],\n  [\"76\u00f78=9, 4\", \"85\u00f74=21, 1\"],\n  [\"19\u00f78=2, 3\", \"34\u00f72=17, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a daily \"two-digit \u00f7 one-digit\" division drill sheet: a\n# centered date heading followed by a 5-column table whose populated rows\n# hold \"<dividend>\u00f7<divisor>=<quotient>, <remainder>\" answers. This commit\n# bumps the date by one day and swaps in a freshly generated set of answers\n# for every populated cell.\n$d = $word.ActiveDocument\n\n# Ordered (old text, new text) pairs. Every old string occurs exactly once\n# in the document, so a non-wildcard Find/Replace-All per pair swaps exactly\n# one run of text while leaving every other property (font/size/etc) intact.\n$replacements = @(\n    @(\"2024-06-02 Sunday\", \"2024-06-03 Monday\"),\n    @(\"20\u00f77=2, 6\", \"93\u00f75=18, 3\"),\n    @(\"95\u00f78=11, 7\", \"73\u00f77=10, 3\"),\n    @(\"46\u00f74=11, 2\", \"59\u00f78=7, 3\"),\n    @(\"81\u00f77=11, 4\", \"78\u00f78=9, 6\"),\n    @(\"68\u00f76=11, 2\", \"72\u00f78=9, 0\"),\n    @(\"87\u00f79=9, 6\", \"15\u00f78=1, 7\"),\n    @(\"48\u00f72=24, 0\", \"91\u00f77=13, 0\"),\n    @(\"55\u00f75=11, 0\", \"25\u00f72=12, 1\"),\n    @(\"69\u00f72=34, 1\", \"23\u00f72=11, 1\"),\n    @(\"33\u00f74=8, 1\", \"75\u00f72=37, 1\"),\n    @(\"13\u00f72=6, 1\", \"12\u00f72=6, 0\"),\n    @(\"53\u00f74=13, 1\", \"24\u00f78=3, 0\"),\n    @(\"25\u00f74=6, 1\", \"59\u00f73=19, 2\"),\n    @(\"35\u00f79=3, 8\", \"36\u00f79=4, 0\"),\n    @(\"54\u00f75=10, 4\", \"37\u00f78=4, 5\"),\n    @(\"18\u00f75=3, 3\", \"90\u00f79=10, 0\"),\n    @(\"30\u00f73=10, 0\", \"96\u00f74=24, 0\"),\n    @(\"70\u00f74=17, 2\", \"96\u00f75=19, 1\"),\n    @(\"20\u00f73=6, 2\", \"43\u00f78=5, 3\"),\n    @(\"39\u00f76=6, 3\", \"55\u00f76=9, 1\"),\n    @(\"78\u00f77=11, 1\", \"93\u00f77=13, 2\"),\n    @(\"84\u00f75=16, 4\", \"76\u00f74=19, 0\"),\n    @(\"13\u00f76=2, 1\", \"93\u00f75=18, 3\"),\n    @(\"76\u00f78=9, 4\", \"85\u00f74=21, 1\"),\n    @(\"19\u00f78=2, 3\", \"34\u00f72=17, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $oldText,  # FindText\n        $false,    # MatchCase\n        $true,     # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap := wdFindContinue\n        $false,    # Format\n        $newText,  # ReplaceWith\n        2          # Replace := wdReplaceAll\n    )\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n$d.Saved = $false\n"}
